$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 226
$ws.Range("B226").Value = 7483281
$ws.Range("F226").Value = "SD Aucas"
$ws.Range("G226").Value = "Delfin SC"
$ws.Range("H226").Value = 0
$ws.Range("I226").Value = 0
$ws.Range("J226").Value = "D"
$ws.Range("K226").Value = 1.909
$ws.Range("L226").Value = 3.25
$ws.Range("M226").Value = 4.2
$ws.Range("N226").Value = 1.909
$ws.Range("O226").Value = 3.5
$ws.Range("P226").Value = 4
$ws.Range("Q226").Value = -0.5
$ws.Range("R226").Value = 1.9
$ws.Range("S226").Value = 1.9
$ws.Range("T226").Value = 2.5
$ws.Range("U226").Value = 1.8
$ws.Range("V226").Value = 2
$ws.Range("W226").Value = -1
$ws.Range("X226").Value = 2.5
$ws.Range("Y226").Value = -1
$ws.Range("Z226").Value = -1
$ws.Range("AA226").Value = 0.8999999999999999
$ws.Range("AB226").Value = -1
$ws.Range("AC226").Value = 1

# Row 227
$ws.Range("B227").Value = 7483189
$ws.Range("F227").Value = "Independiente del Valle"
$ws.Range("G227").Value = "Orense"
$ws.Range("H227").Value = 2
$ws.Range("I227").Value = 2
$ws.Range("J227").Value = "D"
$ws.Range("K227").Value = 1.4
$ws.Range("L227").Value = 4.75
$ws.Range("M227").Value = 7
$ws.Range("N227").Value = 1.4
$ws.Range("O227").Value = 4.5
$ws.Range("P227").Value = 8
$ws.Range("Q227").Value = -1.25
$ws.Range("R227").Value = 1.875
$ws.Range("S227").Value = 1.925
$ws.Range("T227").Value = 2.5
$ws.Range("U227").Value = 1.925
$ws.Range("V227").Value = 1.875
$ws.Range("W227").Value = -1
$ws.Range("X227").Value = 3.5
$ws.Range("Y227").Value = -1
$ws.Range("Z227").Value = -1
$ws.Range("AA227").Value = 0.925
$ws.Range("AB227").Value = 0.925
$ws.Range("AC227").Value = -1

# Row 228
$ws.Range("B228").Value = 7483247
$ws.Range("F228").Value = "Mushuc Runa"
$ws.Range("G228").Value = "Universidad Catolica del Ecuador"
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 2
$ws.Range("J228").Value = "A"
$ws.Range("K228").Value = 3.25
$ws.Range("L228").Value = 3.2
$ws.Range("M228").Value = 2.25
$ws.Range("N228").Value = 3.5
$ws.Range("O228").Value = 3.25
$ws.Range("P228").Value = 2.1
$ws.Range("Q228").Value = 0.5
$ws.Range("R228").Value = 1.775
$ws.Range("S228").Value = 2.025
$ws.Range("T228").Value = 2.5
$ws.Range("U228").Value = 1.9
$ws.Range("V228").Value = 1.9
$ws.Range("W228").Value = -1
$ws.Range("X228").Value = -1
$ws.Range("Y228").Value = 1.1
$ws.Range("Z228").Value = -1
$ws.Range("AA228").Value = 1.025
$ws.Range("AB228").Value = -1
$ws.Range("AC228").Value = 0.8999999999999999

# Row 230
$ws.Range("B230").Value = 7482832
$ws.Range("F230").Value = "Barcelona Guayaquil"
$ws.Range("G230").Value = "Guayaquil City"
$ws.Range("H230").Value = 2
$ws.Range("I230").Value = 1
$ws.Range("J230").Value = "H"
$ws.Range("K230").Value = 1.363
$ws.Range("L230").Value = 5
$ws.Range("M230").Value = 7.5
$ws.Range("N230").Value = 1.444
$ws.Range("O230").Value = 4
$ws.Range("P230").Value = 8
$ws.Range("Q230").Value = -1.25
$ws.Range("R230").Value = 2.05
$ws.Range("S230").Value = 1.75
$ws.Range("T230").Value = 2.5
$ws.Range("U230").Value = 1.95
$ws.Range("V230").Value = 1.85
$ws.Range("W230").Value = 0.444
$ws.Range("X230").Value = -1
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = -0.5
$ws.Range("AA230").Value = 0.375
$ws.Range("AB230").Value = 0.95
$ws.Range("AC230").Value = -1

# Row 232
$ws.Range("B232").Value = 7483188
$ws.Range("F232").Value = "Gualaceo SC"
$ws.Range("G232").Value = "Emelec"
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = 2
$ws.Range("J232").Value = "A"
$ws.Range("K232").Value = 3.6
$ws.Range("L232").Value = 3.3
$ws.Range("M232").Value = 2.05
$ws.Range("N232").Value = 2.6
$ws.Range("O232").Value = 3.25
$ws.Range("P232").Value = 2.75
$ws.Range("Q232").Value = 0
$ws.Range("R232").Value = 1.8
$ws.Range("S232").Value = 2
$ws.Range("T232").Value = 2.5
$ws.Range("U232").Value = 1.975
$ws.Range("V232").Value = 1.825
$ws.Range("W232").Value = -1
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = 1.75
$ws.Range("Z232").Value = -1
$ws.Range("AA232").Value = 1
$ws.Range("AB232").Value = -1
$ws.Range("AC232").Value = 0.825

# Row 233
$ws.Range("B233").Value = 7483306
$ws.Range("F233").Value = "Tecnico Universitario"
$ws.Range("G233").Value = "Club Atletico Libertad"
$ws.Range("H233").Value = 1
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = "D"
$ws.Range("K233").Value = 1.5
$ws.Range("L233").Value = 4.333
$ws.Range("M233").Value = 5.75
$ws.Range("N233").Value = 1.533
$ws.Range("O233").Value = 4.2
$ws.Range("P233").Value = 5.5
$ws.Range("Q233").Value = -1
$ws.Range("R233").Value = 1.925
$ws.Range("S233").Value = 1.875
$ws.Range("T233").Value = 2.25
$ws.Range("U233").Value = 1.8
$ws.Range("V233").Value = 2
$ws.Range("W233").Value = -1
$ws.Range("X233").Value = 3.2
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = -1
$ws.Range("AA233").Value = 0.875
$ws.Range("AB233").Value = -0.5
$ws.Range("AC233").Value = 0.5

# Row 235
$ws.Range("B235").Value = 7528849
$ws.Range("F235").Value = "Guayaquil City"
$ws.Range("G235").Value = "Gualaceo SC"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 2
$ws.Range("J235").Value = "A"
$ws.Range("K235").Value = 1.833
$ws.Range("L235").Value = 3.5
$ws.Range("M235").Value = 3.75
$ws.Range("N235").Value = 2.15
$ws.Range("O235").Value = 3.4
$ws.Range("P235").Value = 3
$ws.Range("Q235").Value = -0.25
$ws.Range("R235").Value = 1.825
$ws.Range("S235").Value = 1.975
$ws.Range("T235").Value = 2.5
$ws.Range("U235").Value = 1.85
$ws.Range("V235").Value = 1.95
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = 2
$ws.Range("Z235").Value = -1
$ws.Range("AA235").Value = 0.9750000000000001
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 0.95

# Row 236
$ws.Range("B236").Value = 7528859
$ws.Range("F236").Value = "Club Atletico Libertad"
$ws.Range("G236").Value = "Cumbaya FC"
$ws.Range("H236").Value = 3
$ws.Range("I236").Value = 1
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 1.727
$ws.Range("L236").Value = 3.5
$ws.Range("M236").Value = 4.333
$ws.Range("N236").Value = 1.4
$ws.Range("O236").Value = 4.2
$ws.Range("P236").Value = 7
$ws.Range("Q236").Value = -1.25
$ws.Range("R236").Value = 2
$ws.Range("S236").Value = 1.8
$ws.Range("T236").Value = 2.5
$ws.Range("U236").Value = 1.95
$ws.Range("V236").Value = 1.85
$ws.Range("W236").Value = 0.3999999999999999
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 1
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = 0.95
$ws.Range("AC236").Value = -1

# Row 238
$ws.Range("B238").Value = 7528857
$ws.Range("F238").Value = "Universidad Catolica del Ecuador"
$ws.Range("G238").Value = "Barcelona Guayaquil"
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 1
$ws.Range("J238").Value = "A"
$ws.Range("K238").Value = 1.533
$ws.Range("L238").Value = 4
$ws.Range("M238").Value = 5.5
$ws.Range("N238").Value = 1.5
$ws.Range("O238").Value = 4.333
$ws.Range("P238").Value = 5.25
$ws.Range("Q238").Value = -1
$ws.Range("R238").Value = 1.8
$ws.Range("S238").Value = 2
$ws.Range("T238").Value = 3
$ws.Range("U238").Value = 1.975
$ws.Range("V238").Value = 1.825
$ws.Range("W238").Value = -1
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = 4.25
$ws.Range("Z238").Value = -1
$ws.Range("AA238").Value = 1
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 0.825

# Row 240
$ws.Range("B240").Value = 7528848
$ws.Range("F240").Value = "Emelec"
$ws.Range("G240").Value = "Deportivo Cuenca"
$ws.Range("H240").Value = 2
$ws.Range("I240").Value = 1
$ws.Range("J240").Value = "H"
$ws.Range("K240").Value = 1.75
$ws.Range("L240").Value = 3.5
$ws.Range("M240").Value = 4.2
$ws.Range("N240").Value = 2.4
$ws.Range("O240").Value = 3.1
$ws.Range("P240").Value = 2.75
$ws.Range("Q240").Value = -0.25
$ws.Range("R240").Value = 2.05
$ws.Range("S240").Value = 1.75
$ws.Range("T240").Value = 2.25
$ws.Range("U240").Value = 1.8
$ws.Range("V240").Value = 2
$ws.Range("W240").Value = 1.4
$ws.Range("X240").Value = -1
$ws.Range("Y240").Value = -1
$ws.Range("Z240").Value = 1.05
$ws.Range("AA240").Value = -1
$ws.Range("AB240").Value = 0.8
$ws.Range("AC240").Value = -1
